$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new "season record" header columns, copying the existing
# header style (bold, bordered, centered) from the last current header cell.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerStyleSource = $ws.Range("AC1")
$headerStyleSource.Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill every data row (2 through 46) with the team's season record.
$wins = 77
$losses = 85
$ties = 0

$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
